$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple per-cell value updates (price / volume refresh) ---
# D-column cells are text-formatted ('@') before assignment so that
# numeric-looking strings (e.g. '556.28') are preserved as text,
# matching the workbook's original inlineStr cell type.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.115.11"
$ws.Range("E2").Value = "  -3.09%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.291.83"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.28"
$ws.Range("E5").Value = "  -3.92%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.32"
$ws.Range("E6").Value = "  -8.14%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.291.58"
$ws.Range("E8").Value = "  -3.95%  "
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.90"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  -5.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.407"
$ws.Range("E12").Value = "  -2.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.852.66"
$ws.Range("E14").Value = "  -0.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.78"
$ws.Range("E15").Value = "  -5.90%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.284.60"
$ws.Range("E16").Value = "  -4.10%  "
$ws.Range("E17").Value = "  -5.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.169.69"
$ws.Range("E18").Value = "  -3.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.08"
$ws.Range("E19").Value = "  -7.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.77"
$ws.Range("E20").Value = "  -5.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.53"
$ws.Range("E21").Value = "  -4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "371.99"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.31"
$ws.Range("E24").Value = "  -4.45%  "
$ws.Range("E25").Value = "  -7.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.416.19"
$ws.Range("E26").Value = "  -4.05%  "
$ws.Range("E27").Value = "  -9.47%  "
$ws.Range("E28").Value = "  -3.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -7.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  -5.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.43"
$ws.Range("E33").Value = "  -5.94%  "
$ws.Range("E34").Value = "  -3.20%  "
$ws.Range("E35").Value = "  -7.72%  "
$ws.Range("E36").Value = "  -8.97%  "
$ws.Range("E37").Value = "  -1.39%  "
$ws.Range("E38").Value = "  -5.47%  "
$ws.Range("E39").Value = "  -5.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.321.00"
$ws.Range("E40").Value = "  -4.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.55"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("E45").Value = "  -4.02%  "
$ws.Range("E46").Value = "  -7.88%  "
$ws.Range("E47").Value = "  -7.13%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.316.50"
$ws.Range("E49").Value = "  -9.27%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.34"
$ws.Range("E50").Value = "  -7.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.47"
$ws.Range("E51").Value = "  -6.14%  "

# --- Row 41 / Row 42 swap: EnergySwap and Hedera swapped ranking order ---
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.83"
$ws.Range("E41").Value = "  -16.95%  "

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0721"
$ws.Range("E42").Value = "  -8.13%  "
